$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.323.42'
$ws.Range("E2").Value = '  -2.36%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.877.69'
$ws.Range("E3").Value = '  -2.31%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '564.25'
$ws.Range("E5").Value = '  -4.90%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.30'
$ws.Range("E6").Value = '  -3.64%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.503'
$ws.Range("E8").Value = '  -0.90%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.876.22'
$ws.Range("E9").Value = '  -2.31%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.80'
$ws.Range("E10").Value = '  -6.67%  '

$ws.Range("E11").Value = '  -3.63%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.430'
$ws.Range("E12").Value = '  -2.97%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000234'
$ws.Range("E13").Value = '  -2.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '31.69'
$ws.Range("E14").Value = '  -3.53%  '

$ws.Range("E15").Value = '  -0.84%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.365.86'
$ws.Range("E16").Value = '  -1.98%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.361.36'
$ws.Range("E17").Value = '  -2.14%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.51'
$ws.Range("E18").Value = '  -2.93%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.888.59'
$ws.Range("E19").Value = '  -3.45%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '429.41'
$ws.Range("E20").Value = '  -2.87%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.02'
$ws.Range("E21").Value = '  -3.30%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.649'
$ws.Range("E22").Value = '  -2.68%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.78'
$ws.Range("E23").Value = '  -3.70%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.85'
$ws.Range("E24").Value = '  -3.09%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.65'
$ws.Range("E25").Value = '  -0.76%  '

$ws.Range("E26").Value = '  +0.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.87'
$ws.Range("E27").Value = '  -11.35%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.99'
$ws.Range("E28").Value = '  -7.19%  '

$ws.Range("E29").Value = '  +0.24%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.94'
$ws.Range("E30").Value = '  -4.07%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.49'
$ws.Range("E31").Value = '  -4.76%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.04'
$ws.Range("E32").Value = '  -9.15%  '

$ws.Range("E33").Value = '  +0.12%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.105'
$ws.Range("E34").Value = '  -3.32%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.39'
$ws.Range("E35").Value = '  -4.14%  '

$ws.Range("E36").Value = '  -3.68%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.35'
$ws.Range("E37").Value = '  -4.50%  '

$ws.Range("E38").Value = '  -1.70%  '

$ws.Range("B39").Value = 'dogwifhat'
$ws.Range("C39").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.81'
$ws.Range("E39").Value = '  -11.56%  '

$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.92'
$ws.Range("E40").Value = '  -5.87%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.18'
$ws.Range("E41").Value = '  -3.78%  '

$ws.Range("E42").Value = '  -3.99%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.13'
$ws.Range("E43").Value = '  -1.58%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.265'
$ws.Range("E44").Value = '  -5.59%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.677.87'
$ws.Range("E45").Value = '  -0.77%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '133.36'
$ws.Range("E46").Value = '  -1.43%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0332'
$ws.Range("E47").Value = '  -1.98%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '336.12'
$ws.Range("E49").Value = '  -7.63%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.102'
$ws.Range("E50").Value = '  -2.32%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.28'
$ws.Range("E51").Value = '  -7.07%  '
